$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 148 values ---
$ws.Range("C148").Value = 0.54
$ws.Range("E148").Value = 0.71
$ws.Range("F148").Value = 2.03
$ws.Range("H148").Value = 1.2
$ws.Range("I148").Value = 2.19

# --- Append new rows 149-170 ---
# Column A holds date-like text; force text format so Excel does not
# auto-convert these to date serials, then restore Normal style so the
# cell does not carry a leftover numeric/date format.
$dates = @(
    "04-08-2021", "05-08-2021", "06-08-2021", "09-08-2021", "10-08-2021", "11-08-2021", "12-08-2021", "13-08-2021", "16-08-2021", "17-08-2021", "18-08-2021", "19-08-2021", "20-08-2021", "23-08-2021", "24-08-2021", "25-08-2021", "26-08-2021", "27-08-2021", "30-08-2021", "31-08-2021", "01-09-2021", "02-09-2021"
)

$startRow = 149
$ws.Range("A149:A170").NumberFormat = "@"
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
}
$ws.Range("A149:A170").Style = "Normal"

# Row 149
$ws.Range("B149").Value = 0.06
$ws.Range("C149").Value = 1.11
$ws.Range("D149").Value = 0.28
$ws.Range("E149").Value = 0.65
$ws.Range("F149").Value = 1.6
$ws.Range("G149").Value = 2
$ws.Range("H149").Value = 0.74
$ws.Range("I149").Value = 2.22

# Row 150
$ws.Range("B150").Value = 0.06
$ws.Range("C150").Value = 0.43
$ws.Range("D150").Value = 0.32
$ws.Range("E150").Value = 0.94
$ws.Range("G150").Value = 2.02
$ws.Range("H150").Value = 2
$ws.Range("I150").Value = 1.06

# Row 151
$ws.Range("B151").Value = 0.06
$ws.Range("C151").Value = 0.97
$ws.Range("D151").Value = 0.24
$ws.Range("E151").Value = 0.77
$ws.Range("F151").Value = 2.75
$ws.Range("G151").Value = 2.04
$ws.Range("H151").Value = 2.1
$ws.Range("I151").Value = 1.44

# Row 152
$ws.Range("B152").Value = 0.06
$ws.Range("C152").Value = 1.08
$ws.Range("D152").Value = 0.23
$ws.Range("E152").Value = 1.11
$ws.Range("F152").Value = 2.16
$ws.Range("G152").Value = 1.75
$ws.Range("H152").Value = 1.2
$ws.Range("I152").Value = 2.45

# Row 153
$ws.Range("C153").Value = 0.67
$ws.Range("D153").Value = 0.31
$ws.Range("E153").Value = 1.12
$ws.Range("F153").Value = 3.43
$ws.Range("G153").Value = 2.15
$ws.Range("H153").Value = 1.57
$ws.Range("I153").Value = 1.52

# Row 154
$ws.Range("B154").Value = 0.06
$ws.Range("C154").Value = 0.97
$ws.Range("D154").Value = 0.31
$ws.Range("E154").Value = 0.69
$ws.Range("F154").Value = 2.09
$ws.Range("G154").Value = 1.81
$ws.Range("H154").Value = 1.5
$ws.Range("I154").Value = 1.46

# Row 155
$ws.Range("B155").Value = 0.06
$ws.Range("C155").Value = 0.51
$ws.Range("D155").Value = 0.25
$ws.Range("E155").Value = 0.57
$ws.Range("F155").Value = 2.61
$ws.Range("G155").Value = 2.23
$ws.Range("H155").Value = 1.13
$ws.Range("I155").Value = 2.26

# Row 156
$ws.Range("B156").Value = 0.06
$ws.Range("C156").Value = 1.52
$ws.Range("D156").Value = 0.29
$ws.Range("E156").Value = 0.82
$ws.Range("G156").Value = 2.24
$ws.Range("H156").Value = 2.29
$ws.Range("I156").Value = 3

# Row 157
$ws.Range("B157").Value = 0.06
$ws.Range("C157").Value = 1.49
$ws.Range("D157").Value = 0.27
$ws.Range("E157").Value = 1.17
$ws.Range("F157").Value = 3.42
$ws.Range("G157").Value = 1.33
$ws.Range("H157").Value = 1.32
$ws.Range("I157").Value = 2.59

# Row 158
$ws.Range("B158").Value = 0.06
$ws.Range("C158").Value = 1.35
$ws.Range("D158").Value = 0.2
$ws.Range("E158").Value = 1.34
$ws.Range("F158").Value = 1.7
$ws.Range("G158").Value = 1.95
$ws.Range("H158").Value = 1.07
$ws.Range("I158").Value = 2.91

# Row 159
$ws.Range("B159").Value = 0.06
$ws.Range("C159").Value = 0.82
$ws.Range("D159").Value = 0.28
$ws.Range("E159").Value = 0.6
$ws.Range("F159").Value = 2.52
$ws.Range("G159").Value = 2.12
$ws.Range("H159").Value = 1.79
$ws.Range("I159").Value = 1.53

# Row 160
$ws.Range("B160").Value = 0.06
$ws.Range("C160").Value = 0.55
$ws.Range("D160").Value = 0.19
$ws.Range("E160").Value = 1.02
$ws.Range("F160").Value = 2.5
$ws.Range("G160").Value = 1.86
$ws.Range("H160").Value = 1.88
$ws.Range("I160").Value = 2.07

# Row 161
$ws.Range("B161").Value = 0.06
$ws.Range("C161").Value = 0.75
$ws.Range("D161").Value = 0.28
$ws.Range("E161").Value = 0.73
$ws.Range("F161").Value = 2.7
$ws.Range("G161").Value = 1.66
$ws.Range("H161").Value = 1.21
$ws.Range("I161").Value = 1.45

# Row 162
$ws.Range("B162").Value = 0.06
$ws.Range("C162").Value = 1.61
$ws.Range("D162").Value = 0.21
$ws.Range("E162").Value = 1.04
$ws.Range("F162").Value = 2.68
$ws.Range("G162").Value = 1.76
$ws.Range("H162").Value = 1.23
$ws.Range("I162").Value = 2.58

# Row 163
$ws.Range("B163").Value = 0.06
$ws.Range("C163").Value = 0.46
$ws.Range("D163").Value = 0.2
$ws.Range("E163").Value = 0.56
$ws.Range("F163").Value = 2.47
$ws.Range("G163").Value = 1.9
$ws.Range("H163").Value = 0.83
$ws.Range("I163").Value = 2.94

# Row 164
$ws.Range("B164").Value = 0.06
$ws.Range("C164").Value = 1.54
$ws.Range("D164").Value = 0.28
$ws.Range("E164").Value = 0.48
$ws.Range("F164").Value = 1.93
$ws.Range("G164").Value = 2.02
$ws.Range("H164").Value = 2.79
$ws.Range("I164").Value = 2.45

# Row 165
$ws.Range("B165").Value = 0.06
$ws.Range("C165").Value = 0.54
$ws.Range("D165").Value = 0.27
$ws.Range("E165").Value = 0.78
$ws.Range("F165").Value = 2.09
$ws.Range("G165").Value = 2.06
$ws.Range("H165").Value = 1.65
$ws.Range("I165").Value = 2.68

# Row 166
$ws.Range("B166").Value = 0.06
$ws.Range("C166").Value = 0.64
$ws.Range("D166").Value = 0.25
$ws.Range("E166").Value = 0.66
$ws.Range("F166").Value = 1.75
$ws.Range("G166").Value = 2.77
$ws.Range("H166").Value = 1.79
$ws.Range("I166").Value = 1.64

# Row 167
$ws.Range("B167").Value = 0.06
$ws.Range("C167").Value = 1.17
$ws.Range("D167").Value = 0.33
$ws.Range("E167").Value = 0.9
$ws.Range("F167").Value = 2.11
$ws.Range("G167").Value = 2
$ws.Range("H167").Value = 1.2
$ws.Range("I167").Value = 2.42

# Row 168
$ws.Range("C168").Value = 0.72
$ws.Range("D168").Value = 0.31
$ws.Range("E168").Value = 0.58
$ws.Range("F168").Value = 2.8
$ws.Range("G168").Value = 2.12
$ws.Range("H168").Value = 1.41
$ws.Range("I168").Value = 1.85

# Row 169
$ws.Range("C169").Value = 1.24
$ws.Range("D169").Value = 0.34
$ws.Range("E169").Value = 0.66
$ws.Range("F169").Value = 2.36
$ws.Range("G169").Value = 2.37
$ws.Range("H169").Value = 1.34
$ws.Range("I169").Value = 2.97

# Row 170
$ws.Range("C170").Value = 0.77
$ws.Range("D170").Value = 0.45
$ws.Range("E170").Value = 0.47
$ws.Range("F170").Value = 2.36
$ws.Range("G170").Value = 2.11
$ws.Range("H170").Value = 0.92
$ws.Range("I170").Value = 2.61

